# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" and "Good Drivers" tables on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table -----------------------------------------------
# Row 4: Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.155.1
$ws.Cells.Item(4, 3).Value = 33872          # Critical Minutes
$ws.Cells.Item(4, 4).Value = 83.5           # Good Roaming Calculation (%)

# Row 5: Realtek 8821CE Wireless LAN 802.11ac PCI-E NIC - 2024.10.139.3
$ws.Cells.Item(5, 3).Value = 405            # Critical Minutes

# Row 6: Realtek 8821CE Wireless LAN 802.11ac PCI-E NIC - 2024.10.138.0
# no longer reported this week -- remove the row entirely (shifts the
# Totals row, and everything below it, up by one).
$ws.Rows.Item(6).Delete()

# Totals row (now row 6 after the delete above)
$ws.Cells.Item(6, 2).Value = 27
$ws.Cells.Item(6, 3).Value = 34279

# --- Good Drivers table ------------------------------------------------
# The "Good Drivers" header (row 12) and the column-header row (row 13)
# already sit in the right place after the row delete above, so only the
# data rows (14-27) need their contents refreshed with this week's numbers
# (new driver versions, resorted by vintage, new vintages filled in).

$goodDrivers = @(
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4', 445055, 99.9, '2024-11-10'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.123.322', 16989, 99.9, '2024-06-30'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.152.0', 1021705, 100, '2024-04-15'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.149.0', 81427, 100, '2023-12-20'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.144.0', 17672, 100, '2023-07-10'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.143.0', 326032, 100, '2023-06-05'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.141.0', 48191, 100, '2023-04-17'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.128.0', 82442, 99.9, '2022-08-29'),
    @('Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.124.0', 11789, 99.9, '2022-07-03'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9', 77849, 99.9, '2021-08-18'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1', 34244, 100, '2021-04-27'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2', 59673, 100, '2020-08-05'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6', 113652, 100, '2020-01-06'),
    @('Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1', 56018, 100, '2019-12-14')
)

$startRow = 14
for ($i = 0; $i -lt $goodDrivers.Count; $i++) {
    $row = $startRow + $i
    $entry = $goodDrivers[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
}
